$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Test on Mura" row (row 23) is now complete: record time spent and notes,
# and mark it with the same "done" fill/wrap style used by the other
# finished rows above it (copy formatting only, from B22).
$ws.Range("B22").Copy()
$ws.Range("B23").PasteSpecial(-4122)

$ws.Range("D23").Value = 2.75
$ws.Range("E23").Value = "Fixed bugs, tested different browsers and devices"

# Move the active selection to where the author left off editing.
$ws.Range("D24").Select()
